# Auto-update draw results: append the 2025-10-31 Pick 4 row (row 45).
# Every column in this sheet is stored as text (t="str"), including the
# date and the numeric-looking "Phase" value, so we force a text
# NumberFormat before writing the values that would otherwise be
# auto-coerced to a date serial / number, then clear the formatting
# again so the new row ends up with the same (default) cell style as
# every other row in the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "2025-10-31"

$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "251031"

$ws.Range("B45").Value = "Pick 4"
$ws.Range("D45").Value = "8-0-7-2"
$ws.Range("E45").Value = "2025-10-31T21:39:16.701+04:00"

# Drop the temporary "@" number format so the new cells keep the
# workbook's default (unstyled) look, matching the rest of the table.
$ws.Range("A45:E45").ClearFormats()

# Dismiss the "number stored as text" warning on the newly-added rows,
# extending the ignored-error range over the new row (mirrors what
# Excel does when the flagged cells already carry dismissed warnings).
$ws.Range("A1:E45").Errors.Item(3).Ignore = $true
